$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mine card (地雷, row 11): reworded effect text to be less ambiguous.
$ws.Range("D11").Value = "有牌进入本牌所在槽位时：本牌所在槽位和对位槽位中所有牌点数减1，然后消灭本牌。"

# Update the active cell selection to D12 (matches the saved view state).
$ws.Range("D12").Select()
